$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 28
$ws.Range("H28").Value = 496.90475
$ws.Range("I28").Value = 555.4706
$ws.Range("J28").Value = 248
$ws.Range("K28").Value = 555.4706
$ws.Range("L28").Value = 248
$ws.Range("M28").Value = -70.47059999999999
$ws.Range("N28").Value = -1218

# Row 43
$ws.Range("H43").Value = 8721.23
$ws.Range("I43").Value = 1379.8
$ws.Range("J43").Value = 13309.625
$ws.Range("K43").Value = 1379.8
$ws.Range("L43").Value = 13309.625
$ws.Range("M43").Value = -1310.8
$ws.Range("N43").Value = -13447.625

# Row 125
$ws.Range("H125").Value = 827.3684
$ws.Range("I125").Value = 706.3077
$ws.Range("J125").Value = 1089.6666
$ws.Range("K125").Value = 6356.7693
$ws.Range("L125").Value = 9806.999400000001
$ws.Range("M125").Value = -3896.7693
$ws.Range("N125").Value = -14726.9994

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 649.2632
$ws.Range("I45").Value = 613.0833
$ws.Range("K45").Value = 613.0833
$ws.Range("M45").Value = -236.0833

# Row 61
$ws.Range("H61").Value = 503581.34
$ws.Range("I61").Value = 371945.94
$ws.Range("J61").Value = 776978
$ws.Range("K61").Value = 371945.94
$ws.Range("L61").Value = 776978
$ws.Range("M61").Value = -371733.94
$ws.Range("N61").Value = -777402

# Row 136
$ws.Range("H136").Value = 503581.34
$ws.Range("I136").Value = 371945.94
$ws.Range("J136").Value = 776978
$ws.Range("K136").Value = 1115837.82
$ws.Range("L136").Value = 2330934
$ws.Range("M136").Value = -1113287.82
$ws.Range("N136").Value = -2336034

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 1276.5454
$ws.Range("I20").Value = 1088.8462
$ws.Range("J20").Value = 1547.6666
$ws.Range("K20").Value = 1088.8462
$ws.Range("L20").Value = 1547.6666
$ws.Range("M20").Value = -841.8462
$ws.Range("N20").Value = -2041.6666

# Row 99
$ws.Range("H99").Value = 2142.9062
$ws.Range("I99").Value = 2411.7778
$ws.Range("K99").Value = 2411.7778
$ws.Range("M99").Value = -913.7777999999998

# Row 107
$ws.Range("H107").Value = 1230.069
$ws.Range("I107").Value = 1269.7084
$ws.Range("J107").Value = 1039.8
$ws.Range("K107").Value = 1269.7084
$ws.Range("L107").Value = 1039.8
$ws.Range("M107").Value = 650.2916
$ws.Range("N107").Value = -4879.8

$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 3596.32
$ws.Range("I58").Value = 5232.5454
$ws.Range("J58").Value = 2310.7144
$ws.Range("K58").Value = 5232.5454
$ws.Range("L58").Value = 2310.7144
$ws.Range("M58").Value = -5029.5454
$ws.Range("N58").Value = -2716.7144

# Row 105
$ws.Range("H105").Value = 888.9286
$ws.Range("I105").Value = 821.3043
$ws.Range("J105").Value = 1200
$ws.Range("K105").Value = 821.3043
$ws.Range("L105").Value = 1200
$ws.Range("M105").Value = 925.6957
$ws.Range("N105").Value = -4694

# Row 132
$ws.Range("H132").Value = 12501959
$ws.Range("I132").Value = 25001280
$ws.Range("J132").Value = 2639.75
$ws.Range("K132").Value = 75003840
$ws.Range("L132").Value = 7919.25
$ws.Range("M132").Value = -75001310
$ws.Range("N132").Value = -12979.25

# Row 134
$ws.Range("H134").Value = 16130289
$ws.Range("I134").Value = 20000834
$ws.Range("J134").Value = 3016.6667
$ws.Range("K134").Value = 60002502
$ws.Range("L134").Value = 9050.000100000001
$ws.Range("M134").Value = -59999967
$ws.Range("N134").Value = -14120.0001

# Row 136
$ws.Range("H136").Value = 3596.32
$ws.Range("I136").Value = 5232.5454
$ws.Range("J136").Value = 2310.7144
$ws.Range("K136").Value = 15697.6362
$ws.Range("L136").Value = 6932.1432
$ws.Range("M136").Value = -13147.6362
$ws.Range("N136").Value = -12032.1432

$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 147.41667
$ws.Range("I2").Value = 289.75
$ws.Range("J2").Value = 76.25
$ws.Range("K2").Value = 1738.5
$ws.Range("L2").Value = 457.5
$ws.Range("M2").Value = -1625.5
$ws.Range("N2").Value = -683.5

# Row 55
$ws.Range("H55").Value = 2758.5293
$ws.Range("J55").Value = 2905.9375
$ws.Range("L55").Value = 8717.8125
$ws.Range("N55").Value = -9071.8125

# Row 82
$ws.Range("H82").Value = 7487.5
$ws.Range("I82").Value = 900
$ws.Range("J82").Value = 8428.571
$ws.Range("K82").Value = 2700
$ws.Range("L82").Value = 25285.713
$ws.Range("M82").Value = -2294
$ws.Range("N82").Value = -26097.713

# Row 85
$ws.Range("H85").Value = 7487.5
$ws.Range("I85").Value = 900
$ws.Range("J85").Value = 8428.571
$ws.Range("K85").Value = 2700
$ws.Range("L85").Value = 25285.713
$ws.Range("M85").Value = -1296
$ws.Range("N85").Value = -28093.713

# Row 92
$ws.Range("H92").Value = 666.6667
$ws.Range("I92").Value = 400
$ws.Range("J92").Value = 800
$ws.Range("K92").Value = 1200
$ws.Range("L92").Value = 2400
$ws.Range("M92").Value = 48
$ws.Range("N92").Value = -4896

# Row 113
$ws.Range("H113").Value = 644.0769
$ws.Range("I113").Value = 606
$ws.Range("J113").Value = 696
$ws.Range("K113").Value = 1818
$ws.Range("L113").Value = 2088
$ws.Range("M113").Value = 352
$ws.Range("N113").Value = -6428

# Row 122
$ws.Range("H122").Value = 794.1111
$ws.Range("I122").Value = 329.4
$ws.Range("J122").Value = 1375
$ws.Range("K122").Value = 2964.6
$ws.Range("L122").Value = 12375
$ws.Range("M122").Value = -514.5999999999999
$ws.Range("N122").Value = -17275

# Row 131
$ws.Range("H131").Value = 1177.8518
$ws.Range("J131").Value = 1257.35
$ws.Range("L131").Value = 3772.05
$ws.Range("N131").Value = -13852.05

$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 2121.7334
$ws.Range("I102").Value = 1677.1
$ws.Range("J102").Value = 3011
$ws.Range("K102").Value = 1677.1
$ws.Range("L102").Value = 3011
$ws.Range("M102").Value = -55.09999999999991
$ws.Range("N102").Value = -6255

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 832.2222
$ws.Range("I22").Value = 648.3333
$ws.Range("J22").Value = 1200
$ws.Range("K22").Value = 648.3333
$ws.Range("L22").Value = 1200
$ws.Range("M22").Value = -353.3333
$ws.Range("N22").Value = -1790

# Row 27
$ws.Range("H27").Value = 832.2222
$ws.Range("I27").Value = 648.3333
$ws.Range("J27").Value = 1200
$ws.Range("K27").Value = 648.3333
$ws.Range("L27").Value = 1200
$ws.Range("M27").Value = -541.3333
$ws.Range("N27").Value = -1414

# Row 46
$ws.Range("H46").Value = 1017.95654
$ws.Range("I46").Value = 856.125
$ws.Range("J46").Value = 1387.8572
$ws.Range("K46").Value = 856.125
$ws.Range("L46").Value = 1387.8572
$ws.Range("M46").Value = -668.125
$ws.Range("N46").Value = -1763.8572

# Row 134
$ws.Range("H134").Value = 46896
$ws.Range("J134").Value = 46896
$ws.Range("L134").Value = 46896
$ws.Range("N134").Value = -57036
